# Replace the "Authors" column (E) values for rows 2-13 with cleaned-up
# versions of the author lists (whitespace/formatting cleanup of the
# underlying scraped author data, per the commit message).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '[Michela%Casanova%xref no email%1,           Elena%Pagani Bagliacca%xref no email%1,           Matteo%Silva%xref no email%1,           Carlo%Patriarca%xref no email%1,           Laura%Veneroni%xref no email%1,           Carlo Alfredo%Clerici%xref no email%1,           Filippo%Spreafico%xref no email%1,           Roberto%Luksch%xref no email%1,           Monica%Terenziani%xref no email%1,           Cristina%Meazza%xref no email%1,           Marta%Podda%xref no email%1,           Veronica%Biassoni%xref no email%1,           Elisabetta%Schiavello%xref no email%1,           Stefano%Chiaravalli%xref no email%1,           Nadia%Puma%xref no email%1,           Luca%Bergamaschi%xref no email%1,           Giovanna%Gattuso%xref no email%1,           Giovanna%Sironi%xref no email%1,           Maura%Massimino%xref no email%1,           Andrea%Ferrari%xref no email%1]'
$ws.Range("E3").Value = '[Marco%Colizzi%NULL%1,                  Elena%Sironi%NULL%2,                  Elena%Sironi%NULL%0,                  Federico%Antonini%NULL%1,                  Marco Luigi%Ciceri%NULL%1,                  Chiara%Bovo%NULL%2,                  Chiara%Bovo%NULL%0,                  Leonardo%Zoccante%NULL%1]'
$ws.Range("E4").Value = '[Wen Yan%Jiao%NULL%1,                  Lin Na%Wang%NULL%1,                  Juan%Liu%NULL%1,                  Shuan Feng%Fang%NULL%1,                  Fu Yong%Jiao%NULL%1,                  Massimo%Pettoello-Mantovani%NULL%1,                  Eli%Somekh%NULL%1]'
$ws.Range("E5").Value = '[Shuang-Jiang%Zhou%NULL%0,                  Li-Gang%Zhang%NULL%3,                  Lei-Lei%Wang%NULL%3,                  Zhao-Chang%Guo%NULL%3,                  Jing-Qi%Wang%NULL%3,                  Jin-Cheng%Chen%NULL%3,                  Mei%Liu%NULL%3,                  Xi%Chen%NULL%0,                  Jing-Xu%Chen%chenjx1110@163.com%3]'
$ws.Range("E6").Value = '[Xi%Liu%NULL%1,                  Wen-Tao%Luo%NULL%1,                  Ying%Li%NULL%0,                  Chun-Na%Li%NULL%1,                  Zhong-Si%Hong%NULL%1,                  Hui-Li%Chen%NULL%1,                  Fei%Xiao%NULL%1,                  Jin-Yu%Xia%xiajinyu@mail.sysu.edu.cn%1]'
$ws.Range("E7").Value = '[Benjamin%Oosterhoff%Benjamin.oosterhoff@montana.edu%2,                  Cara A.%Palmer%NULL%2,                  Jenna%Wilson%NULL%2,                  Natalie%Shook%NULL%2]'
$ws.Range("E8").Value = '[İsmail%Seçer%ismailsecer84@gmail.com%2,                  Sümeyye%Ulaş%NULL%4,                  Sümeyye%Ulaş%NULL%0]'
$ws.Range("E9").Value = '[Fangyuan%Tian%NULL%1,                  Hongxia%Li%lihx@xust.edu.cn%1,                  Shuicheng%Tian%tiansc@xust.edu.cn%1,                  Jie%Yang%NULL%0,                  Jiang%Shao%NULL%1,                  Chenning%Tian%NULL%1]'
$ws.Range("E10").Value = '[Abrams%D.%coreGivesNoEmail%1,                 Antonovsky%A.%coreGivesNoEmail%1,                 Bowlby%J.%coreGivesNoEmail%2,                 Bowlby%J.%coreGivesNoEmail%0,                 Durkheim%\u00c9.%coreGivesNoEmail%1,                 Erikson%E. H.%coreGivesNoEmail%1,                 Haidt%J.%coreGivesNoEmail%1,                 Janoff-Bulman%R.%coreGivesNoEmail%1,                 Jerzy%Trzebi\u0144ski%coreGivesNoEmail%1,                 Jolanta%Zuzanna Czarnecka%coreGivesNoEmail%1,                 Maciej%Caba\u0144ski%coreGivesNoEmail%1,                 Seligman%M. E. P.%coreGivesNoEmail%2,                 Seligman%M. E. P.%coreGivesNoEmail%0,                 Spielberger%C. D.%coreGivesNoEmail%1,                 Trzebi\u0144ski%J.%coreGivesNoEmail%1,                 Zi\u0119ba%M.%coreGivesNoEmail%2,                 Zi\u0119ba%M.%coreGivesNoEmail%0]'
$ws.Range("E11").Value = '[Huan%Yang%xref no email%1,           Peng%Bin%xref no email%1,           Alex Jingwei%He%xref no email%1]'
$ws.Range("E12").Value = '[Xinyan%Xie%xref no email%1,           Qi%Xue%xref no email%1,           Yu%Zhou%xref no email%1,           Kaiheng%Zhu%xref no email%1,           Qi%Liu%xref no email%1,           Jiajia%Zhang%xref no email%1,           Ranran%Song%xref no email%1]'
$ws.Range("E13").Value = '[Shuang-Jiang%Zhou%NULL%0,                  Li-Gang%Zhang%NULL%0,                  Lei-Lei%Wang%NULL%0,                  Zhao-Chang%Guo%NULL%0,                  Jing-Qi%Wang%NULL%0,                  Jin-Cheng%Chen%NULL%0,                  Mei%Liu%NULL%0,                  Xi%Chen%NULL%0,                  Jing-Xu%Chen%chenjx1110@163.com%0]'

Write-Output "done"
